# Update Top50_DataComp.xlsx data: extend M2 coverage by one additional
# month for the rows whose underlying source refreshed since the last
# pull (bumping the M2 observation count in column C, or the M2 first
# date in column E, and always pushing the M2 last date in column F
# forward by one month).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (ColumnLetter => NewValue) for the M2_Len / M2_1stDate bump.
$lenUpdates = @{
    2  = 330
    14 = 419
    16 = 450
    17 = 365
    21 = 305
    23 = 255
    24 = 211
    31 = 389
    33 = 450
    34 = 197
    36 = 498
    41 = 389
    42 = 234
    43 = 468
    44 = 401
    48 = 353
    52 = 342
}

$firstDateUpdates = @{
    5  = "1981-11-01"
    7  = "1981-10-01"
    9  = "1981-10-01"
    30 = "1981-10-01"
}

# Every touched row also gets its M2_LastDate (column F) pushed one
# month later.
$lastDateUpdates = @{
    2  = "2023-06-01"
    5  = "2023-06-01"
    7  = "2023-05-01"
    9  = "2023-05-01"
    14 = "2023-05-01"
    16 = "2023-05-01"
    17 = "2023-05-01"
    21 = "2023-05-01"
    23 = "2023-04-01"
    24 = "2023-06-01"
    30 = "2023-05-01"
    31 = "2023-05-01"
    33 = "2023-06-01"
    34 = "2023-05-01"
    36 = "2023-06-01"
    41 = "2023-05-01"
    42 = "2023-05-01"
    43 = "2023-05-01"
    44 = "2023-05-01"
    48 = "2023-05-01"
    52 = "2023-06-01"
}

foreach ($row in $lenUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $lenUpdates[$row]
}

foreach ($row in $firstDateUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $firstDateUpdates[$row]
}

foreach ($row in $lastDateUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $lastDateUpdates[$row]
}
